$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 14:22"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1064737
$ws.Range("C4").Value = 543
$ws.Range("E4").Value = 855656
$ws.Range("F4").Value = 18851
$ws.Range("G4").Value = 14
$ws.Range("H4").Value = 61670

# --- Alemania (row 9) ---
$ws.Range("B9").Value = 161552
$ws.Range("C9").Value = 13
$ws.Range("E9").Value = 31585

# --- Paises Bajos (row 17) ---
$ws.Range("B17").Value = 39316
$ws.Range("C17").Value = 514
$ws.Range("E17").Value = 34271
$ws.Range("G17").Value = 84
$ws.Range("H17").Value = 4795

# --- Suecia (row 24) ---
$ws.Range("B24").Value = 21092
$ws.Range("C24").Value = 790
$ws.Range("E24").Value = 17501
$ws.Range("F24").Value = 531
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = 2586

# --- Dinamarca (row 41) ---
$ws.Range("D41").Value = 6546
$ws.Range("E41").Value = 2160
$ws.Range("F41").Value = 62
$ws.Range("G41").Value = 9
$ws.Range("H41").Value = 452

# --- Croacia / Ghana / Armenia reshuffle (rows 68-70) ---
# Croacia is inserted into the country list right before Ghana, so it now
# occupies row 68 with freshly updated figures; the old Ghana and Armenia
# rows shift down one position (their figures unchanged) to rows 69 and 70.
$ws.Range("A68").Value = "Croacia"
$ws.Range("B68").Value = 2076
$ws.Range("C68").Value = 14
$ws.Range("D68").Value = 1348
$ws.Range("E68").Value = 659
$ws.Range("F68").Value = 19
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 69

$ws.Range("A69").Value = "Ghana"
$ws.Range("B69").Value = 2074
$ws.Range("C69").Value = 403
$ws.Range("D69").Value = 212
$ws.Range("E69").Value = 1845
$ws.Range("F69").Value = 4
$ws.Range("G69").Value = 1
$ws.Range("H69").Value = 17

$ws.Range("A70").Value = "Armenia"
$ws.Range("B70").Value = 2066
$ws.Range("C70").Value = 134
$ws.Range("D70").Value = 929
$ws.Range("E70").Value = 1105
$ws.Range("F70").Value = 10
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 32

# --- Maldivas (row 128) ---
$ws.Range("B128").Value = 301
$ws.Range("C128").Value = 23
$ws.Range("E128").Value = 283
